$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D:E) for the new reporting quarters (2018-09-30, 2018-12-31).
# This shifts the existing quarterly data (previously in D:K) right into F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# Fix up number formats on the freshly inserted D:E columns by copying the
# formatting from column F (which now holds what used to be column D), but
# only for the rows that actually carry quarterly data (skip section-title
# rows, which have no D:K cells at all).
$ws.Range("F7:F35").Copy() | Out-Null
$ws.Range("D7:E35").PasteSpecial(-4122) | Out-Null
$ws.Range("F38:F77").Copy() | Out-Null
$ws.Range("D38:E77").PasteSpecial(-4122) | Out-Null
$ws.Range("F80:F102").Copy() | Out-Null
$ws.Range("D80:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the two new quarter columns with their reported figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 2957000
$ws.Range("E8").Value = 2807000
$ws.Range("D9").Value = 288000
$ws.Range("E9").Value = 228000
$ws.Range("D10").Value = 2669000
$ws.Range("E10").Value = 2579000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 80000
$ws.Range("E15").Value = 78000
$ws.Range("D17").Value = 1747000
$ws.Range("E17").Value = 1588000
$ws.Range("D18").Value = 1210000
$ws.Range("E18").Value = 1219000
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 1290000
$ws.Range("E21").Value = 1297000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 1210000
$ws.Range("E23").Value = 1219000
$ws.Range("D24").Value = 275000
$ws.Range("E24").Value = 296000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 935000
$ws.Range("E26").Value = 923000
$ws.Range("D27").Value = 885000
$ws.Range("E27").Value = 885000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 885000
$ws.Range("E33").Value = 885000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 885000
$ws.Range("E35").Value = 885000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 27938000
$ws.Range("E41").Value = 21830000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 22204000
$ws.Range("E43").Value = 23209000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 20758000
$ws.Range("E45").Value = 4063000
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 220540000
$ws.Range("E47").Value = 217998000
$ws.Range("D48").Value = 1769000
$ws.Range("E48").Value = 1683000
$ws.Range("D49").Value = 1379000
$ws.Range("E49").Value = 1227000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 296482000
$ws.Range("E54").Value = 272102000
$ws.Range("D57").Value = 34557000
$ws.Range("E57").Value = 29373000
$ws.Range("D58").Value = 231423000
$ws.Range("E58").Value = 213408000
$ws.Range("D59").Value = 2954000
$ws.Range("E59").Value = 2697000
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 6878000
$ws.Range("E61").Value = 5790000
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 275812000
$ws.Range("E66").Value = 251268000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 2793000
$ws.Range("E70").Value = 2793000
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 17329000
$ws.Range("E72").Value = 16615000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 17877000
$ws.Range("E76").Value = 18041000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 885000
$ws.Range("E81").Value = 885000
$ws.Range("D83").Value = 80000
$ws.Range("E83").Value = 78000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 7801000
$ws.Range("E89").Value = 1275000
$ws.Range("D91").Value = -170000
$ws.Range("E91").Value = -147000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -14535000
$ws.Range("E94").Value = -6366000
$ws.Range("D96").Value = -208000
$ws.Range("E96").Value = -225000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 17867000
$ws.Range("E100").Value = 13267000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 11133000
$ws.Range("E102").Value = 8176000

# A handful of previously-reported figures were restated/corrected for older quarters.
$ws.Range("F45").Value = 5621000
$ws.Range("G45").Value = 8389000
$ws.Range("H45").Value = 8543000
$ws.Range("I45").Value = 8686000
$ws.Range("J45").Value = 10895000
$ws.Range("F47").Value = 214799000
$ws.Range("G47").Value = 198833000
$ws.Range("H47").Value = 194534000
$ws.Range("I47").Value = 186344000
$ws.Range("J47").Value = 177109000
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("H89").Value = -1188000
$ws.Range("I89").Value = 1942000
$ws.Range("I91").Value = -103000
$ws.Range("J91").Value = -84000
$ws.Range("H102").Value = 1341000
$ws.Range("I102").Value = 3301000

